$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.750.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.636.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.08%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.76%  "

$ws.Range("E9").Value = "  +2.13%  "

$ws.Range("E10").Value = "  +0.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0860"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.867.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.640.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.25%  "

$ws.Range("E14").Value = "  -0.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.555"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.722.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0720"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("E20").Value = "  -1.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.30%  "

$ws.Range("E23").Value = "  +1.41%  "

$ws.Range("E24").Value = "  +3.09%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("E26").Value = "  -0.68%  "

$ws.Range("E27").Value = "  -0.56%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  -0.09%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.16%  "

$ws.Range("E31").Value = "  -0.43%  "

$ws.Range("E32").Value = "  -0.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.482.63"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.72%  "

$ws.Range("E34").Value = "  -0.81%  "

$ws.Range("E35").Value = "  -1.29%  "

$ws.Range("E36").Value = "  -1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.958"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.28%  "

$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.559"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.33%  "

$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "67.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.63%  "

$ws.Range("B43").Value = "mCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.02%  "

$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.777.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.41%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.22%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.12%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.75%  "
